# Gantt Chard update
# - PERIODS value G3: 7 -> 9
# - Row 9 (Airfoil Selection): Actual Duration F9 "-" -> 2, Percent Complete G9 0.5 -> 1
# - Row 11 (Wing Design and Lift Analysis): Actual Start E11 "-" -> 7, Actual Duration F11 "-" -> 2, Percent Complete G11 0 -> 1
# - Row 12 (Tail Design and Stability Analysis): Actual Start E12 "-" -> 8, Percent Complete G12 0 -> 0.3
# - Selection cursor moved to E15

$wb = $excel.ActiveWorkbook
$ws = $wb.ActiveSheet

# PERIODS
$ws.Range("G3").Value = 9

# Airfoil Selection (row 9)
$ws.Range("F9").Value = 2
$ws.Range("G9").Value = 1

# Wing Design and Lift Analysis (row 11)
$ws.Range("E11").Value = 7
$ws.Range("F11").Value = 2
$ws.Range("G11").Value = 1

# Tail Design and Stability Analysis (row 12)
$ws.Range("E12").Value = 8
$ws.Range("G12").Value = 0.3

# Move selection cursor
$ws.Range("E15").Select() | Out-Null
